# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" column-header suffixes to the respective
# input-file-version suffixes ("_FV2310" / "_FV2404"), turns the header
# row + data range into a native Excel Table ("Table1"), and freezes the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1. Rename header cells: "..._old" -> "..._FV2310", "..._new" -> "..._FV2404" ---
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $name = $cell.Value()
    if ($name -ne $null) {
        if ($name -like "*_old") {
            $cell.Value = ($name -replace "_old$", "_FV2310")
        } elseif ($name -like "*_new") {
            $cell.Value = ($name -replace "_new$", "_FV2404")
        }
    }
}

# --- 2. Turn the range into a Table, preserving the existing header formatting ---
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))

# Stash the current header formatting on a scratch row far outside the used
# range, so it can be restored after Table creation without the engine
# baking a header dxf / default table style into the workbook.
$scratchRow = $lastRow + 111
$helperRange = $ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow, $lastCol))
$headerRange.Copy()
$helperRange.PasteSpecial(-4122)   # xlPasteFormats

$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

$helperRange.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats
$helperRange.ClearFormats()
$helperRange.ClearContents()
$excel.CutCopyMode = $false

$tbl.TableStyle = ""

# --- 3. Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
